# Update the documentation with for 'Tabella comuni subentrati'
#
# 1) Remove the "SQL Statement" worksheet entirely (it only held the raw
#    SQL query text, which is no longer published). This also drops the
#    now-unused SQL-text shared string and re-numbers the header on
#    "SQL Results" (G1) back onto the plain "DATASUBENTRO" string.
# 2) Append new "comuni subentrati" rows (12..14) to the results table,
#    matching the existing formatting of rows 10/11 (style index 6 for
#    columns A/D/E, style index 2 for columns B/C/F, style index 3 for
#    column G).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$wsResults = $wb.Worksheets.Item("SQL Results")

# --- Remove the "SQL Statement" sheet -------------------------------------
$wsStatement = $wb.Worksheets.Item("SQL Statement")
$wsStatement.Delete()

# --- Append the new data rows ------------------------------------------
# Rows 10/11 already exist (untouched) using a style pattern -- A/D/E use
# style index 6 ("text"/quote-prefix), B/C/F use style index 2
# ("bordered"), G uses style index 3 ("bordered date"). Row 10 is reused
# below purely as a formatting template (PasteSpecial xlPasteFormats) so
# that no new cell styles get introduced into the workbook.
$templateRow = 10

$rows = @(
    @{ Row = 12; A = "098019"; B = "C816"; C = "CODOGNO";  D = 98; E = "03"; F = "LO"; G = 42951.736805555556 },
    @{ Row = 13; A = "054053"; B = "L216"; C = "TORGIANO"; D = 54; E = "10"; F = "PG"; G = 42970.736805555556 },
    @{ Row = 14; A = "010046"; B = "H183"; C = "RAPALLO";  D = 10; E = "07"; F = "GE"; G = 42972.736805555556 }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Fill in the values first. Columns A/C/F/B/E are written in this
    # particular order to reproduce the shared-string table ordering of
    # the reference workbook. Leading-zero / zero-padded codes (columns A
    # and E) need a leading apostrophe so they are kept as text instead
    # of being coerced into numbers (e.g. "098019" -> 98019, "03" -> 3).
    $wsResults.Cells.Item($row, 1).Value = "'" + $r.A
    $wsResults.Cells.Item($row, 3).Value = $r.C
    $wsResults.Cells.Item($row, 6).Value = $r.F
    $wsResults.Cells.Item($row, 2).Value = $r.B
    $wsResults.Cells.Item($row, 5).Value = "'" + $r.E
    $wsResults.Cells.Item($row, 4).Value = $r.D
    $wsResults.Cells.Item($row, 7).Value = $r.G

    # ... then copy the formatting from the template row onto the new row
    # (per cell, so we don't blow the whole 16384-column row out). Doing
    # this after the values are set preserves the quote-prefixed "text"
    # styling (style index 6) on numeric cells such as column D, exactly
    # like the reference workbook.
    for ($col = 1; $col -le 7; $col++) {
        $wsResults.Cells.Item($templateRow, $col).Copy()
        $wsResults.Cells.Item($row, $col).PasteSpecial(-4122)
    }
}

$excel.CutCopyMode = $false

# --- Update the selection bookkeeping --------------------------------------
$wsResults.Range("E17").Select()
